# Update attendee/view counts (column F) on the "展览" and "全部类型" sheets.
# Both sheets share the same underlying rows, with F40 differing by 1
# between the two sheets in the updated data.

$wb = $excel.ActiveWorkbook

# Common updates shared by both "展览" and "全部类型" sheets (row -> new F value)
$commonUpdates = @{
    2  = 23
    6  = 112
    8  = 452
    9  = 46
    10 = 21
    11 = 565
    12 = 28
    15 = 366
    18 = 3
    19 = 48
    20 = 48
    21 = 93
    22 = 907
    23 = 1393
    24 = 295
    25 = 322
    27 = 72
    28 = 153
    31 = 211
    32 = 241
    33 = 269
    34 = 1612
    37 = 155
    38 = 577
    41 = 420
    42 = 189
    43 = 906
}

function Update-SheetColumnF {
    param(
        $SheetName,
        $Row40Value
    )

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($row in $commonUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $commonUpdates[$row]
    }

    # Row 40 differs between the two sheets
    $ws.Cells.Item(40, 6).Value = $Row40Value
}

Update-SheetColumnF "展览" 3591
Update-SheetColumnF "全部类型" 3592
